# AddShift.xlsx — refresh the test-run results for this pass.
# Five previously-executed rows ("Y") are reset back to "N" (not yet run),
# and the selection is moved over to the J5:L5 block (Revise column area)
# to match where the author was working next.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reset the "Execute" flag on the rows whose test results were re-run.
$ws.Range("A6").Value = "N"
$ws.Range("A10").Value = "N"
$ws.Range("A18").Value = "N"
$ws.Range("A26").Value = "N"
$ws.Range("A27").Value = "N"

# Move the selection/view over towards the J:L columns (Result/Revise area).
$ws.Range("J5:L5").Select() | Out-Null
